$wb = $excel.ActiveWorkbook
$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("D2").Value = "2016-03-23 03:52:53"
$zhcn.Range("E2").Value = "2016-03-23 03:52:41"
$dede.Range("E2").Value = "2016-03-23 03:52:53"
